# Two TM5 variables (emibvoc, emivoc) have been identified by Tommi Bergman
# as now available in TM5, so they are removed from the list of ignored
# CMIP6-requested variables.
#
# In the worksheet these correspond to the AERmon rows for:
#   emibvoc - total emission rate of biogenic nmvoc
#   emivoc  - total emission rate of nmvoc
# which sit at rows 14 and 15. Deleting the rows shifts every row below
# them up by two (AERmonZ block, LPJ-GUESS block, trailing blank rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("14:15").Delete()

# Author's commit left the active selection on A14 (the row that now
# contains the next variable, "hcho", after the deleted rows).
$ws.Range("A14").Select()
